# Auto-generated edit script applying the cryptos.xlsx price-refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.294.69"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.610.78"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  -0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "212.94"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.42%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "18.42"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.96%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0815"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "1.833.59"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "1.594.24"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("E14").Value = "  +0.67%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.516"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "26.292.39"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "61.89"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  -0.13%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "203.48"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +8.36%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "144.25"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -2.80%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.26"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("E30").Value = "  +3.92%  "
$ws.Range("E31").Value = "  -0.13%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.21"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.11%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.95"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("E34").Value = "  +3.34%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").Value = "1.159.05"
$ws.Range("E36").Value = "  +4.67%  "
$ws.Range("E37").Value = "  +9.73%  "
$ws.Range("B38").Value = "PaxDollar"
$ws.Range("C38").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.797"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("E40").Value = "  -0.54%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.503"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.72%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.785"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.35%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.97%  "
$ws.Range("D44").Value = "1.746.60"
$ws.Range("E44").Value = "  +0.28%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "91.77"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.71%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.53"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.05%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "54.42"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₇0966"
$ws.Range("E49").Value = "  -10.67%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.406"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("E51").Value = "  -0.27%  "
